$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.520.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.241.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  +1.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.37%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.81'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0803'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("E13").Value = '  +0.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.836'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.219.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '44.208.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("E18").Value = '  -0.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '65.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("E23").Value = '  +0.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.68%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  +3.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0797'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.27%  '

$ws.Range("E35").Value = '  +3.57%  '

$ws.Range("E36").Value = '  +0.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.91'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.40%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0303'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.67%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.62%  '

$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.794.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.192'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.68%  '

$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '79.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.51%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '70.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '99.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.83%  '
